$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3085878193809322
$ws.Cells.Item(2, 3).Value = 0.05600425874996517
$ws.Cells.Item(2, 5).Value = 0.7412275368263153
$ws.Cells.Item(2, 6).Value = 2.1401641376887
$ws.Cells.Item(2, 7).Value = 0.3757734516185138
$ws.Cells.Item(2, 8).Value = 0.5512079301942805
$ws.Cells.Item(2, 10).Value = 0.03500250578657216
$ws.Cells.Item(2, 11).Value = 0.2883637197925282
$ws.Cells.Item(2, 14).Value = 1.195404981206032
$ws.Cells.Item(2, 15).Value = 1.795432379671979

$ws.Cells.Item(3, 2).Value = 0.2729701417155752
$ws.Cells.Item(3, 3).Value = 0.05234576605603536
$ws.Cells.Item(3, 5).Value = 0.7121726695951338
$ws.Cells.Item(3, 6).Value = 2.110706947619377
$ws.Cells.Item(3, 7).Value = 0.3789724841568187
$ws.Cells.Item(3, 8).Value = 0.5559883586556182
$ws.Cells.Item(3, 10).Value = 0.03533973316615402
$ws.Cells.Item(3, 11).Value = 0.2518252571454127
$ws.Cells.Item(3, 14).Value = 1.205710661661431
$ws.Cells.Item(3, 15).Value = 1.812005848290596

$ws.Cells.Item(4, 2).Value = 0.2510848548860736
$ws.Cells.Item(4, 3).Value = 0.05008712320821473
$ws.Cells.Item(4, 5).Value = 0.6947055943195437
$ws.Cells.Item(4, 6).Value = 2.093884989288568
$ws.Cells.Item(4, 7).Value = 0.3812152619849272
$ws.Cells.Item(4, 8).Value = 0.559159155452349
$ws.Cells.Item(4, 10).Value = 0.03556788842769798
$ws.Cells.Item(4, 11).Value = 0.2293227894513592
$ws.Cells.Item(4, 14).Value = 1.212459214350723
$ws.Cells.Item(4, 15).Value = 1.823258438441982

$ws.Cells.Item(5, 2).Value = 0.2421629972761536
$ws.Cells.Item(5, 3).Value = 0.04916367737625649
$ws.Cells.Item(5, 5).Value = 0.6876812727684438
$ws.Cells.Item(5, 6).Value = 2.087347678862642
$ws.Cells.Item(5, 7).Value = 0.3821991548976769
$ws.Cells.Item(5, 8).Value = 0.560510543242188
$ws.Cells.Item(5, 10).Value = 0.0356661634232367
$ws.Cells.Item(5, 11).Value = 0.2201363969329577
$ws.Cells.Item(5, 14).Value = 1.21531520039948
$ws.Cells.Item(5, 15).Value = 1.82811445182324

$ws.Cells.Item(6, 2).Value = 0.2406813395037943
$ws.Cells.Item(6, 3).Value = 0.04901015880128057
$ws.Cells.Item(6, 5).Value = 0.6865205439177373
$ws.Cells.Item(6, 6).Value = 2.086281347893021
$ws.Cells.Item(6, 7).Value = 0.382366750913917
$ws.Cells.Item(6, 8).Value = 0.5607385200570008
$ws.Cells.Item(6, 10).Value = 0.03568280186612327
$ws.Cells.Item(6, 11).Value = 0.2186100264622581
$ws.Cells.Item(6, 14).Value = 1.215795832866633
$ws.Cells.Item(6, 15).Value = 1.828937120875764

$ws.Cells.Item(7, 2).Value = 0.2509645447092055
$ws.Cells.Item(7, 3).Value = 0.05007468147744021
$ws.Cells.Item(7, 5).Value = 0.6946104827594439
$ws.Cells.Item(7, 6).Value = 2.093795538565203
$ws.Cells.Item(7, 7).Value = 0.3812282480502986
$ws.Cells.Item(7, 8).Value = 0.5591771407629906
$ws.Cells.Item(7, 10).Value = 0.03556919234596734
$ws.Cells.Item(7, 11).Value = 0.229198964415275
$ws.Cells.Item(7, 14).Value = 1.212497302305721
$ws.Cells.Item(7, 15).Value = 1.823322833369119

$ws.Cells.Item(8, 2).Value = 0.2963104727025723
$ws.Cells.Item(8, 3).Value = 0.05474540451390908
$ws.Cells.Item(8, 5).Value = 0.731132003633121
$ws.Cells.Item(8, 6).Value = 2.129744583998303
$ws.Cells.Item(8, 7).Value = 0.3768186007994032
$ws.Cells.Item(8, 8).Value = 0.5528073365306341
$ws.Cells.Item(8, 10).Value = 0.03511440026133172
$ws.Cells.Item(8, 11).Value = 0.2757796410004687
$ws.Cells.Item(8, 14).Value = 1.198871083953051
$ws.Cells.Item(8, 15).Value = 1.800923392460206

$ws.Cells.Item(9, 2).Value = 0.385086857421868
$ws.Cells.Item(9, 3).Value = 0.06380450004247962
$ws.Cells.Item(9, 5).Value = 0.8057158946062373
$ws.Cells.Item(9, 6).Value = 2.210296528264564
$ws.Cells.Item(9, 7).Value = 0.3703863872033395
$ws.Cells.Item(9, 8).Value = 0.5421846171184583
$ws.Cells.Item(9, 10).Value = 0.03439017504906339
$ws.Cells.Item(9, 11).Value = 0.3665662045767988
$ws.Cells.Item(9, 14).Value = 1.17548566752194
$ws.Cells.Item(9, 15).Value = 1.765547667567802

$ws.Cells.Item(10, 2).Value = 0.4501997104115105
$ws.Cells.Item(10, 3).Value = 0.0703964307267313
$ws.Cells.Item(10, 5).Value = 0.8623378670615693
$ws.Cells.Item(10, 6).Value = 2.275645480384952
$ws.Cells.Item(10, 7).Value = 0.3670180179778129
$ws.Cells.Item(10, 8).Value = 0.5355179679821802
$ws.Cells.Item(10, 10).Value = 0.03396060355742136
$ws.Cells.Item(10, 11).Value = 0.4329052416027253
$ws.Cells.Item(10, 14).Value = 1.160332976982389
$ws.Cells.Item(10, 15).Value = 1.744781794817641

$ws.Cells.Item(11, 2).Value = 0.4797928833435208
$ws.Cells.Item(11, 3).Value = 0.07338088227047024
$ws.Cells.Item(11, 5).Value = 0.8884971670285751
$ws.Cells.Item(11, 6).Value = 2.306722228280279
$ws.Cells.Item(11, 7).Value = 0.365781927139146
$ws.Cells.Item(11, 8).Value = 0.5327319747677848
$ws.Cells.Item(11, 10).Value = 0.03378751375275257
$ws.Cells.Item(11, 11).Value = 0.4630017952918024
$ws.Cells.Item(11, 14).Value = 1.153879127946176
$ws.Cells.Item(11, 15).Value = 1.736472154081511

$ws.Cells.Item(12, 2).Value = 0.490994668479118
$ws.Cells.Item(12, 3).Value = 0.07450890826351042
$ws.Cells.Item(12, 5).Value = 0.8984609711751119
$ws.Cells.Item(12, 6).Value = 2.318684697481814
$ws.Cells.Item(12, 7).Value = 0.3653565649342525
$ws.Cells.Item(12, 8).Value = 0.5317124522613881
$ws.Cells.Item(12, 10).Value = 0.03372518544068726
$ws.Cells.Item(12, 11).Value = 0.4743863719922672
$ws.Cells.Item(12, 14).Value = 1.151498316443359
$ws.Cells.Item(12, 15).Value = 1.733489218644209

$ws.Cells.Item(13, 2).Value = 0.4885823740456772
$ws.Cells.Item(13, 3).Value = 0.07426606315780759
$ws.Cells.Item(13, 5).Value = 0.8963125136550474
$ws.Cells.Item(13, 6).Value = 2.316099714738527
$ws.Cells.Item(13, 7).Value = 0.3654462723766159
$ws.Cells.Item(13, 8).Value = 0.5319304470521971
$ws.Cells.Item(13, 10).Value = 0.03373846577014206
$ws.Cells.Item(13, 11).Value = 0.4719350583094695
$ws.Cells.Item(13, 14).Value = 1.152008259973854
$ws.Cells.Item(13, 15).Value = 1.734124360311739

$ws.Cells.Item(14, 2).Value = 0.48071455486334
$ws.Cells.Item(14, 3).Value = 0.07347372856764878
$ws.Cells.Item(14, 5).Value = 0.8893157354075925
$ws.Cells.Item(14, 6).Value = 2.30770249005252
$ws.Cells.Item(14, 7).Value = 0.3657460754716908
$ws.Cells.Item(14, 8).Value = 0.5326473871247828
$ws.Cells.Item(14, 10).Value = 0.03378232145106885
$ws.Cells.Item(14, 11).Value = 0.4639386620792152
$ws.Cells.Item(14, 14).Value = 1.153681992292142
$ws.Cells.Item(14, 15).Value = 1.736223462965654

$ws.Cells.Item(15, 2).Value = 0.4758946828943635
$ws.Cells.Item(15, 3).Value = 0.07298812240361485
$ws.Cells.Item(15, 5).Value = 0.8850375423392762
$ws.Cells.Item(15, 6).Value = 2.302584272359212
$ws.Cells.Item(15, 7).Value = 0.3659352804695786
$ws.Cells.Item(15, 8).Value = 0.5330911531444684
$ws.Cells.Item(15, 10).Value = 0.03380960347872453
$ws.Cells.Item(15, 11).Value = 0.4590390134669917
$ws.Cells.Item(15, 14).Value = 1.154715420948243
$ws.Cells.Item(15, 15).Value = 1.737530556307817

$ws.Cells.Item(16, 2).Value = 0.4482651303242164
$ws.Cells.Item(16, 3).Value = 0.07020109689668175
$ws.Cells.Item(16, 5).Value = 0.8606363852588004
$ws.Cells.Item(16, 6).Value = 2.273641722258247
$ws.Cells.Item(16, 7).Value = 0.3671047704250157
$ws.Cells.Item(16, 8).Value = 0.5357050031947352
$ws.Cells.Item(16, 10).Value = 0.03397236518963354
$ws.Cells.Item(16, 11).Value = 0.4309366691553578
$ws.Cells.Item(16, 14).Value = 1.160763585824526
$ws.Cells.Item(16, 15).Value = 1.745347748647745

$ws.Cells.Item(17, 2).Value = 0.4313079362780456
$ws.Cells.Item(17, 3).Value = 0.06848764590574774
$ws.Cells.Item(17, 5).Value = 0.8457700049601584
$ws.Cells.Item(17, 6).Value = 2.25623222450767
$ws.Cells.Item(17, 7).Value = 0.367898167439769
$ws.Cells.Item(17, 8).Value = 0.5373716992622946
$ws.Cells.Item(17, 10).Value = 0.03407793689011562
$ws.Cells.Item(17, 11).Value = 0.4136754932501958
$ws.Cells.Item(17, 14).Value = 1.164586400516527
$ws.Cells.Item(17, 15).Value = 1.750434709215583

$ws.Cells.Item(18, 2).Value = 0.4215521036814778
$ws.Cells.Item(18, 3).Value = 0.0675007771430387
$ws.Cells.Item(18, 5).Value = 0.8372570592048874
$ws.Cells.Item(18, 6).Value = 2.246345711046075
$ws.Cells.Item(18, 7).Value = 0.3683823804019397
$ws.Cells.Item(18, 8).Value = 0.5383535593403082
$ws.Cells.Item(18, 10).Value = 0.03414075981222808
$ws.Cells.Item(18, 11).Value = 0.4037396965356379
$ws.Cells.Item(18, 14).Value = 1.166826522482609
$ws.Cells.Item(18, 15).Value = 1.753467585001957

$ws.Cells.Item(19, 2).Value = 0.4182485355683809
$ws.Cells.Item(19, 3).Value = 0.06716641288289793
$ws.Cells.Item(19, 5).Value = 0.8343812122169254
$ws.Cells.Item(19, 6).Value = 2.243020103359299
$ws.Cells.Item(19, 7).Value = 0.3685511097941756
$ws.Cells.Item(19, 8).Value = 0.5386899885283967
$ws.Cells.Item(19, 10).Value = 0.03416239119320252
$ws.Cells.Item(19, 11).Value = 0.4003743160814679
$ws.Cells.Item(19, 14).Value = 1.167592090549086
$ws.Cells.Item(19, 15).Value = 1.754512832265917

$ws.Cells.Item(20, 2).Value = 0.4331133221273546
$ws.Cells.Item(20, 3).Value = 0.06867018459104202
$ws.Cells.Item(20, 5).Value = 0.8473486434521789
$ws.Cells.Item(20, 6).Value = 2.258072351763332
$ws.Cells.Item(20, 7).Value = 0.3678108234145299
$ws.Cells.Item(20, 8).Value = 0.5371918733946117
$ws.Cells.Item(20, 10).Value = 0.03406648112088639
$ws.Cells.Item(20, 11).Value = 0.4155137685946784
$ws.Cells.Item(20, 14).Value = 1.164175177073879
$ws.Cells.Item(20, 15).Value = 1.749882118520262

$ws.Cells.Item(21, 2).Value = 0.4830256513866971
$ws.Cells.Item(21, 3).Value = 0.07370651457796384
$ws.Cells.Item(21, 5).Value = 0.8913692882119193
$ws.Cells.Item(21, 6).Value = 2.310163680148577
$ws.Cells.Item(21, 7).Value = 0.3656568556647031
$ws.Cells.Item(21, 8).Value = 0.5324358417456452
$ws.Cells.Item(21, 10).Value = 0.0337693526108449
$ws.Cells.Item(21, 11).Value = 0.4662877372632579
$ws.Cells.Item(21, 14).Value = 1.153188663702927
$ws.Cells.Item(21, 15).Value = 1.735602459514311

$ws.Cells.Item(22, 2).Value = 0.5156196823012351
$ws.Cells.Item(22, 3).Value = 0.07698565117813416
$ws.Cells.Item(22, 5).Value = 0.9204765291077308
$ws.Cells.Item(22, 6).Value = 2.345341538255951
$ws.Cells.Item(22, 7).Value = 0.3644981678657402
$ws.Cells.Item(22, 8).Value = 0.5295342611150886
$ws.Cells.Item(22, 10).Value = 0.03359391650585586
$ws.Cells.Item(22, 11).Value = 0.4993992718409572
$ws.Cells.Item(22, 14).Value = 1.146376246434123
$ws.Cells.Item(22, 15).Value = 1.727224432598192

$ws.Cells.Item(23, 2).Value = 0.4982262716917489
$ws.Cells.Item(23, 3).Value = 0.07523667217891727
$ws.Cells.Item(23, 5).Value = 0.9049105644993745
$ws.Cells.Item(23, 6).Value = 2.326462646374182
$ws.Cells.Item(23, 7).Value = 0.3650937513660892
$ws.Cells.Item(23, 8).Value = 0.5310639713389662
$ws.Cells.Item(23, 10).Value = 0.03368583175609885
$ws.Cells.Item(23, 11).Value = 0.481733835546521
$ws.Cells.Item(23, 14).Value = 1.149978508769962
$ws.Cells.Item(23, 15).Value = 1.731608513351873

$ws.Cells.Item(24, 2).Value = 0.4322971286798065
$ws.Cells.Item(24, 3).Value = 0.06858766440100794
$ws.Cells.Item(24, 5).Value = 0.8466348352947506
$ws.Cells.Item(24, 6).Value = 2.257240048864162
$ws.Cells.Item(24, 7).Value = 0.3678502241947825
$ws.Cells.Item(24, 8).Value = 0.5372730989866668
$ws.Cells.Item(24, 10).Value = 0.03407165364562736
$ws.Cells.Item(24, 11).Value = 0.4146827220060629
$ws.Cells.Item(24, 14).Value = 1.16436095931742
$ws.Cells.Item(24, 15).Value = 1.750131607398046

$ws.Cells.Item(25, 2).Value = 0.361088452502031
$ws.Cells.Item(25, 3).Value = 0.06136478089071318
$ws.Cells.Item(25, 5).Value = 0.7852197202114297
$ws.Cells.Item(25, 6).Value = 2.187424711902707
$ws.Cells.Item(25, 7).Value = 0.371888606474208
$ws.Cells.Item(25, 8).Value = 0.5448584389979487
$ws.Cells.Item(25, 10).Value = 0.03456811813288851
$ws.Cells.Item(25, 11).Value = 0.3420680156249034
$ws.Cells.Item(25, 14).Value = 1.181455462712606
$ws.Cells.Item(25, 15).Value = 1.774201147620417
